$wb = $excel.ActiveWorkbook

# --- DATA sheet: insert a new "Browser" column before the existing "First Name" column ---
$ws = $wb.Worksheets.Item("DATA")
$ws.Columns.Item(5).Insert()

$ws.Range("E1").Value = "Browser"
$ws.Range("E2").Value = "chrome"
$ws.Range("E3").Value = "chrome"
$ws.Range("E4").Value = "chrome"
$ws.Range("E5").Value = "chrome"

# Fix the "Execute" flag that previously reused the "No" shared string
$ws.Range("B3").Value = "Yes"
$ws.Range("B5").Value = "Yes"

# --- Selections (cursor position) on both sheets ---
$ws1 = $wb.Worksheets.Item("RunManager")
$ws1.Activate() | Out-Null
$ws1.Range("B5").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("E4").Select() | Out-Null
